$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Andrey's score
$ws.Range("B2").Value = 10

# Row 3 becomes Jenivaldo (score stays 6)
$ws.Range("A3").Value = "Jenivaldo"
$ws.Range("B3").Value = 6

# Re-insert Marcelo as row 4
$ws.Range("A4").Value = "Marcelo"
$ws.Range("B4").Value = 8

# Add new approved students
$ws.Range("A5").Value = "Marcos"
$ws.Range("B5").Value = 6

$ws.Range("A6").Value = "Matias"
$ws.Range("B6").Value = 7
